$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 2.82
$ws.Range("B2").Value = 0.026

$ws.Range("B5").Value = 17.61
$ws.Range("C5").Value = 1.9

$ws.Range("B6").Value = 27.82
$ws.Range("C6").Value = 2.9

$ws.Range("B7").Value = 48.58
$ws.Range("C7").Value = 4.9000000000000004

$ws.Range("B3").Select()
